# Auto-generated Excel COM-interop script to update cryptos.xlsx values
# Applies the cell-content changes described in the commit diff (93 cell updates
# across rows 2-51: price (D), volume/% change (E), and two coin name/link swaps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.666.34"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "3.027.62"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'581.55"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "'148.86"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.24%  "
$ws.Range("D9").Value = "3.027.09"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("E10").Value = "  -4.16%  "
$ws.Range("D11").Value = "'5.66"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "'0.0000231"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("D14").Value = "'35.21"
$ws.Range("E14").Value = "  -5.99%  "
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "3.532.00"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "62.653.89"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "'7.02"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "3.029.59"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "'468.21"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "'14.05"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'0.690"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").Value = "'7.38"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'81.04"
$ws.Range("D26").Value = "'12.44"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").Value = "'10.36"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'7.21"
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "'0.108"
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "0.0₃0796"
$ws.Range("E36").Value = "  -7.18%  "
$ws.Range("D37").Value = "'5.78"
$ws.Range("E37").Value = "  -5.00%  "
$ws.Range("D38").Value = "'2.15"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "'50.26"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  -16.43%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.97"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "'421.12"
$ws.Range("E42").Value = "  -6.33%  "
$ws.Range("D43").Value = "'0.113"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "'0.280"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0355"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.792.18"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'38.09"
$ws.Range("E47").Value = "  -10.02%  "
$ws.Range("D48").Value = "'128.76"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.44"
$ws.Range("E50").Value = "  -4.56%  "
$ws.Range("E51").Value = "  -1.77%  "
